$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Periodo Mora" (E), "Valor Mora" (F) and "Salario Basico" (G)
# for the account-statement rows (16-24): reorder periods chronologically
# (2110, 2110, 2111, 2112, 2201, 2202, 2203, 2204, 2205), refresh the
# "Valor Mora" figures and bump "Salario Basico" up to 1,000,000 across
# the board.

$ws.Range("E16").Value = "2110"
$ws.Range("F16").Value = 36341
$ws.Range("G16").Value = 1000000

$ws.Range("E17").Value = "2110"
$ws.Range("F17").Value = 40000
$ws.Range("G17").Value = 1000000

$ws.Range("E18").Value = "2111"
$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 1000000

$ws.Range("E19").Value = "2112"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 1000000

$ws.Range("E20").Value = "2201"
$ws.Range("F20").Value = 40000
$ws.Range("G20").Value = 1000000

$ws.Range("E21").Value = "2202"
$ws.Range("F21").Value = 40000
$ws.Range("G21").Value = 1000000

$ws.Range("E22").Value = "2203"
$ws.Range("F22").Value = 40000
$ws.Range("G22").Value = 1000000

$ws.Range("E23").Value = "2204"
$ws.Range("F23").Value = 40000
$ws.Range("G23").Value = 1000000

$ws.Range("E24").Value = "2205"
$ws.Range("F24").Value = 30666
$ws.Range("G24").Value = 1000000

# Nudge the logo image to the left to re-center it after the data refresh.
# (Compute the new Left precisely from the column A width plus the exact
# target cell offset, rather than doing lossy arithmetic on the already
# rounded .Left getter, so the anchor lands exactly where intended.)
$shp = $ws.Shapes.Item(1)
$colA = $ws.Columns.Item(1).Width
$targetFromColOffEmu = 426600
$shp.Left = $colA + ($targetFromColOffEmu / 12700.0)
